$wb = $excel.ActiveWorkbook

# --- ALC: 54 cell updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2978.6155
$ws.Range("I64").Value = 2869.7222
$ws.Range("J64").Value = 3113.7932
$ws.Range("K64").Value = 2869.7222
$ws.Range("L64").Value = 3113.7932
$ws.Range("M64").Value = -2621.7222
$ws.Range("N64").Value = -3609.7932
$ws.Range("H67").Value = 2978.6155
$ws.Range("I67").Value = 2869.7222
$ws.Range("J67").Value = 3113.7932
$ws.Range("K67").Value = 2869.7222
$ws.Range("L67").Value = 3113.7932
$ws.Range("M67").Value = -2011.7222
$ws.Range("N67").Value = -4829.7932
$ws.Range("H92").Value = 22733410
$ws.Range("I92").Value = 35718788
$ws.Range("J92").Value = 8996.25
$ws.Range("K92").Value = 35718788
$ws.Range("L92").Value = 8996.25
$ws.Range("M92").Value = -35717540
$ws.Range("N92").Value = -11492.25
$ws.Range("H94").Value = 8441.333000000001
$ws.Range("J94").Value = 10528
$ws.Range("L94").Value = 10528
$ws.Range("N94").Value = -11430
$ws.Range("H96").Value = 316.07144
$ws.Range("I96").Value = 258.70587
$ws.Range("K96").Value = 776.11761
$ws.Range("M96").Value = 596.88239
$ws.Range("H100").Value = 55557188
$ws.Range("I100").Value = 1742.8572
$ws.Range("J100").Value = 250001250
$ws.Range("K100").Value = 1742.8572
$ws.Range("L100").Value = 250001250
$ws.Range("M100").Value = -1201.8572
$ws.Range("N100").Value = -250002332
$ws.Range("H137").Value = 2101.55
$ws.Range("I137").Value = 1252.9286
$ws.Range("J137").Value = 4081.6667
$ws.Range("K137").Value = 3758.7858
$ws.Range("L137").Value = 12245.0001
$ws.Range("M137").Value = -1208.7858
$ws.Range("N137").Value = -17345.0001
$ws.Range("H138").Value = 2466.48
$ws.Range("J138").Value = 2659.24
$ws.Range("L138").Value = 7977.719999999999
$ws.Range("N138").Value = -18257.72
$ws.Range("H141").Value = 2802.7334
$ws.Range("I141").Value = 3014.1
$ws.Range("J141").Value = 2380
$ws.Range("K141").Value = 9042.299999999999
$ws.Range("L141").Value = 7140
$ws.Range("M141").Value = -3862.299999999999
$ws.Range("N141").Value = -17500

# --- ARM: 29 cell updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3183.4
$ws.Range("I63").Value = 1393
$ws.Range("J63").Value = 4750
$ws.Range("K63").Value = 1393
$ws.Range("L63").Value = 4750
$ws.Range("M63").Value = -707
$ws.Range("N63").Value = -6122
$ws.Range("H66").Value = 3183.4
$ws.Range("I66").Value = 1393
$ws.Range("J66").Value = 4750
$ws.Range("K66").Value = 6965
$ws.Range("L66").Value = 23750
$ws.Range("M66").Value = -3533
$ws.Range("N66").Value = -30614
$ws.Range("H104").Value = 21891.666
$ws.Range("J104").Value = 21891.666
$ws.Range("L104").Value = 21891.666
$ws.Range("N104").Value = -28879.666
$ws.Range("H111").Value = 21733.334
$ws.Range("J111").Value = 21733.334
$ws.Range("L111").Value = 21733.334
$ws.Range("N111").Value = -29913.334
$ws.Range("H122").Value = 1483.4
$ws.Range("I122").Value = 1370.4445
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4111.333500000001
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1661.333500000001
$ws.Range("N122").Value = -12400

# --- CRP: 26 cell updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1645.2084
$ws.Range("I58").Value = 1496.1765
$ws.Range("J58").Value = 2007.1428
$ws.Range("K58").Value = 1496.1765
$ws.Range("L58").Value = 2007.1428
$ws.Range("M58").Value = -1293.1765
$ws.Range("N58").Value = -2413.1428
$ws.Range("H62").Value = 3635.2942
$ws.Range("I62").Value = 3612.5
$ws.Range("K62").Value = 3612.5
$ws.Range("M62").Value = -2988.5
$ws.Range("H65").Value = 3635.2942
$ws.Range("I65").Value = 3612.5
$ws.Range("K65").Value = 18062.5
$ws.Range("M65").Value = -14942.5
$ws.Range("H135").Value = 34748.094
$ws.Range("J135").Value = 34748.094
$ws.Range("L135").Value = 34748.094
$ws.Range("N135").Value = -44888.094
$ws.Range("H136").Value = 1645.2084
$ws.Range("I136").Value = 1496.1765
$ws.Range("J136").Value = 2007.1428
$ws.Range("K136").Value = 4488.529500000001
$ws.Range("L136").Value = 6021.428400000001
$ws.Range("M136").Value = -1938.529500000001
$ws.Range("N136").Value = -11121.4284

# --- CUL: 28 cell updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 400.66666
$ws.Range("I5").Value = 369.07144
$ws.Range("J5").Value = 511.25
$ws.Range("K5").Value = 1107.21432
$ws.Range("L5").Value = 1533.75
$ws.Range("M5").Value = -995.21432
$ws.Range("N5").Value = -1757.75
$ws.Range("H70").Value = 5992.6816
$ws.Range("I70").Value = 5074.75
$ws.Range("J70").Value = 7094.2
$ws.Range("K70").Value = 15224.25
$ws.Range("L70").Value = 21282.6
$ws.Range("M70").Value = -14909.25
$ws.Range("N70").Value = -21912.6
$ws.Range("H73").Value = 5992.6816
$ws.Range("I73").Value = 5074.75
$ws.Range("J73").Value = 7094.2
$ws.Range("K73").Value = 15224.25
$ws.Range("L73").Value = 21282.6
$ws.Range("M73").Value = -14132.25
$ws.Range("N73").Value = -23466.6
$ws.Range("H135").Value = 400.66666
$ws.Range("I135").Value = 369.07144
$ws.Range("J135").Value = 511.25
$ws.Range("K135").Value = 3321.64296
$ws.Range("L135").Value = 4601.25
$ws.Range("M135").Value = -786.6429600000001
$ws.Range("N135").Value = -9671.25

# --- GSM: 4 cell updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2245.0908
$ws.Range("I122").Value = 2458
$ws.Range("K122").Value = 7374
$ws.Range("M122").Value = -4924

# --- LTW: 11 cell updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 737.5
$ws.Range("I55").Value = 160
$ws.Range("J55").Value = 1150
$ws.Range("K55").Value = 160
$ws.Range("L55").Value = 1150
$ws.Range("M55").Value = 13
$ws.Range("N55").Value = -1496
$ws.Range("H93").Value = 2014.375
$ws.Range("I93").Value = 1815.3334
$ws.Range("K93").Value = 1815.3334
$ws.Range("M93").Value = -567.3334

# --- WVR: 7 cell updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2857.8
$ws.Range("I132").Value = 3141
$ws.Range("J132").Value = 2433
$ws.Range("K132").Value = 9423
$ws.Range("L132").Value = 7299
$ws.Range("M132").Value = -6893
$ws.Range("N132").Value = -12359
